$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H24").Value = 2
$ws.Range("D25").Value = "Em andamento"
$ws.Range("H25").Value = 1
$ws.Range("H26").Value = 1

$ws.Range("D26").Select()
